# Commit message: "nonshared aggregation = composition"
#
# The only real content change in the target diff is in the speaker notes
# of the "Aggregation: Shared vs. Non-shared" slide: a new paragraph is
# appended after the existing "Shared aggregation - many parts can be
# associated with Whole..." bullet, introducing Composition as an even
# stronger, non-shared form of aggregation.
#
# (The rest of the reference diff is just PowerPoint renumbering the
# customXml/item*.xml SharePoint metadata parts on save - not a real
# content edit, so nothing to replicate for it here.)

$p = $ppt.ActivePresentation

# Find the slide titled "Aggregation: Shared vs. Non-shared" (slide 18 in
# this deck) robustly, in case slide ordering ever shifts.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidateSlide = $p.Slides.Item($i)
    if ($candidateSlide.Shapes.HasTitle) {
        $title = $candidateSlide.Shapes.Title.TextFrame.TextRange.Text
        if ($title -eq "Aggregation: Shared vs. Non-shared") {
            $targetSlide = $candidateSlide
        }
    }
}
if ($targetSlide -eq $null) {
    $targetSlide = $p.Slides.Item(18)
}

$notes = $targetSlide.NotesPage

# Locate the notes body placeholder (holds the "Shared aggregation ..."
# text) robustly by name, falling back to the known index.
$target = $null
for ($i = 1; $i -le $notes.Shapes.Count; $i++) {
    $candidate = $notes.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text -like "Shared aggregation*") {
            $target = $candidate
        }
    }
}
if ($target -eq $null) {
    $target = $notes.Shapes.Item(4)
}

$tr = $target.TextFrame.TextRange
$existingText = $tr.Text

# Build the new paragraph's text (runs get merged by the host when a
# notes TextRange is rewritten, so compose the final combined string):
#   "Composition" + " is an even stronger form of "non‑shared" aggregation,"
# Chr(34) = straight double quote, Chr(8209) = non-breaking hyphen (U+2011).
$quote = [char]34
$nbHyphen = [char]8209
$newParagraph = "Composition is an even stronger form of " + $quote + "non" + $nbHyphen + "shared" + $quote + " aggregation,"

# Append as its own paragraph. A bare line-feed (Chr(10)) makes the host
# emit a separate <a:p> for the new paragraph instead of a literal break
# character inside the existing run's text.
$tr.Text = $existingText + [char]10 + $newParagraph
